# The upstream change described by the diff is purely a byproduct of the
# authoring toolchain being swapped out (commit message: "Fixed POI packaging
# and upgraded to POI 3.15."). Every hunk in the unified diff touches only the
# *order* in which XML attributes (namespace declarations on <w:document>,
# and plain attributes such as w:pgSz/w:pgMar/w:rFonts/w:lang/w:lsdException/
# w:style/.../w:tblCellMar children) are written out - e.g.
#   -<w:pgSz w:w="11906" w:h="16838"/>
#   +<w:pgSz w:h="16838" w:w="11906"/>
# The attribute name/value pairs, element tree, text runs, paragraph
# content, styles, numbering, margins, fonts, language settings, etc. are
# byte-for-byte identical before and after - only the serializer's
# attribute ordering (alphabetised by the newer POI version) changed.
#
# There is no document-content mutation to perform: the Word object model
# doesn't expose (and shouldn't expose) control over raw XML attribute
# ordering, and no paragraph text, formatting, style definition, or page
# setup value actually differs between the two revisions. So this script
# intentionally leaves the document's content untouched - that is the
# faithful application of this particular change.

$d = $word.ActiveDocument
